$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column D holds prices formatted as plain text (e.g. '56.419.21',
# '1.00'), some of which look like numbers/dates to Excel's auto-
# detection. Prefix with a literal apostrophe to force text entry,
# then reset the style so no stray number-format / quote-prefix
# styling is left behind on the cell.
function Set-TextValue($cell, $value) {
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '56.419.21'
$ws.Range("E2").Value = '  -3.55%  '

Set-TextValue $ws.Range("D3") '2.974.32'
$ws.Range("E3").Value = '  -5.65%  '

$ws.Range("E4").Value = '  +0.12%  '

Set-TextValue $ws.Range("D5") '495.27'
$ws.Range("E5").Value = '  -6.11%  '

Set-TextValue $ws.Range("D6") '134.74'
$ws.Range("E6").Value = '  -0.41%  '

Set-TextValue $ws.Range("D7") '1.00'
$ws.Range("E7").Value = '  -0.03%  '

Set-TextValue $ws.Range("D8") '0.424'
$ws.Range("E8").Value = '  -4.72%  '

Set-TextValue $ws.Range("D9") '7.17'
$ws.Range("E9").Value = '  -2.00%  '

$ws.Range("E10").Value = '  -3.97%  '

$ws.Range("E11").Value = '  -7.79%  '

$ws.Range("E12").Value = '  -0.66%  '

Set-TextValue $ws.Range("D13") '3.487.79'
$ws.Range("E13").Value = '  -5.41%  '

Set-TextValue $ws.Range("D14") '25.15'
$ws.Range("E14").Value = '  -1.24%  '

Set-TextValue $ws.Range("D15") '56.444.47'
$ws.Range("E15").Value = '  -3.37%  '

Set-TextValue $ws.Range("D16") '2.975.00'
$ws.Range("E16").Value = '  -5.39%  '

$ws.Range("E17").Value = '  -4.24%  '

Set-TextValue $ws.Range("D18") '5.80'
$ws.Range("E18").Value = '  +0.28%  '

Set-TextValue $ws.Range("D19") '12.30'
$ws.Range("E19").Value = '  -5.98%  '

Set-TextValue $ws.Range("D20") '7.75'
$ws.Range("E20").Value = '  -2.37%  '

Set-TextValue $ws.Range("D21") '324.80'
$ws.Range("E21").Value = '  -5.50%  '

$ws.Range("E22").Value = '  -0.06%  '

Set-TextValue $ws.Range("D23") '0.468'
$ws.Range("E23").Value = '  -8.35%  '

Set-TextValue $ws.Range("D24") '61.47'
$ws.Range("E24").Value = '  -8.75%  '

Set-TextValue $ws.Range("D25") '0.998'
$ws.Range("E25").Value = '  -0.11%  '

$ws.Range("E26").Value = '  -5.76%  '

Set-TextValue $ws.Range("D27") '0.0₃0889'
$ws.Range("E27").Value = '  -6.91%  '

$ws.Range("E28").Value = '  -0.18%  '

Set-TextValue $ws.Range("D29") '6.49'
$ws.Range("E29").Value = '  -5.66%  '

Set-TextValue $ws.Range("D30") '6.74'
$ws.Range("E30").Value = '  -2.82%  '

$ws.Range("E31").Value = '  -7.10%  '

$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range("D32") '20.29'
$ws.Range("E32").Value = '  -5.48%  '

$ws.Range("B33").Value = 'Fetch.AI'
$ws.Range("C33").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws.Range("D33") '1.16'
$ws.Range("E33").Value = '  -7.08%  '

Set-TextValue $ws.Range("D34") '151.99'
$ws.Range("E34").Value = '  -4.05%  '

Set-TextValue $ws.Range("D35") '4.44'
$ws.Range("E35").Value = '  -8.63%  '

$ws.Range("E36").Value = '  -7.97%  '

Set-TextValue $ws.Range("D37") '5.58'
$ws.Range("E37").Value = '  -11.03%  '

Set-TextValue $ws.Range("D38") '0.0668'
$ws.Range("E38").Value = '  -2.90%  '

Set-TextValue $ws.Range("D39") '23.09'
$ws.Range("E39").Value = '  -4.11%  '

Set-TextValue $ws.Range("D40") '3.009.78'
$ws.Range("E40").Value = '  -5.35%  '

Set-TextValue $ws.Range("D41") '36.64'
$ws.Range("E41").Value = '  -9.52%  '

Set-TextValue $ws.Range("D42") '1.00'
$ws.Range("E42").Value = '  +0.17%  '

Set-TextValue $ws.Range("D43") '0.638'
$ws.Range("E43").Value = '  -8.14%  '

$ws.Range("E44").Value = '  -9.52%  '

Set-TextValue $ws.Range("D45") '2.220.85'
$ws.Range("E45").Value = '  -2.78%  '

$ws.Range("E46").Value = '  -4.63%  '

Set-TextValue $ws.Range("D47") '3.54'
$ws.Range("E47").Value = '  -9.80%  '

Set-TextValue $ws.Range("D48") '1.93'
$ws.Range("E48").Value = '  +3.26%  '

Set-TextValue $ws.Range("D49") '0.0236'
$ws.Range("E49").Value = '  +0.36%  '

Set-TextValue $ws.Range("D50") '5.72'
$ws.Range("E50").Value = '  -7.10%  '

Set-TextValue $ws.Range("D51") '18.86'
$ws.Range("E51").Value = '  -9.23%  '

